$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "mollie_ABC123"
$ws.Range("B3").Value = 1445758
$ws.Range("C3").Value = "Ben Gortemaker"
$ws.Range("D3").Value = 1
